# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" (summary) sheet,
#    cloning the layout/format of the "2021-Q4" sheet, and fill it with the
#    2022-Q1 per-fund holdings.
# 2) Insert a new top row into "总计" for the 2022-Q1 aggregate figures and
#    renumber the existing index column.
#
# NOTE: worksheet handles returned by Item()/Add() in this COM layer are
# positional, not stable object identities - once the sheet collection is
# mutated (Add/Insert/Delete), any previously-fetched sheet variable can
# silently resolve to a *different* sheet afterwards. To stay safe we
# re-fetch sheets by name with $wb.Worksheets.Item(...) immediately before
# using them whenever the sheet count may have changed in between.

$wb = $excel.ActiveWorkbook

# --- 1) create the new "2022-Q1" sheet, placed right before "总计" -----
$tot = $wb.Worksheets.Item("总计")
$q1new = $wb.Worksheets.Add($tot)
$q1new.Name = "2022-Q1"

# Re-fetch everything by name now that the sheet collection changed size -
# any handle obtained before the Add() call above may now point elsewhere.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1new = $wb.Worksheets.Item("2022-Q1")

# Clone formatting (styles/borders/fonts) from the 2021-Q4 sheet so the new
# sheet matches the established look (bold/bordered header row + index col).
$q4.Range("A1:H10").Copy()
$q1new.Range("A1").PasteSpecial(-4122)

# Header row
$q1new.Range("B1").Value = "基金代码"
$q1new.Range("C1").Value = "基金名称"
$q1new.Range("D1").Value = "基金规模"
$q1new.Range("E1").Value = "股票总仓位"
$q1new.Range("F1").Value = "仓位占比"
$q1new.Range("G1").Value = "持有市值(亿元)"
$q1new.Range("H1").Value = "仓位排名"

# Helper data for the 9 funds held in 2022-Q1.
$rows = @(
    @(0, "013414", "太平智远三个月定期开放股票", "8.69", "86.34", "6.85", "0.5953", 3),
    @(1, "011738", "华安兴安优选一年持有期混合型证券投资基金A", "25.77", "54.03", "1.45", "0.3737", 5),
    @(2, "011739", "华安兴安优选一年持有期混合型证券投资基金C", "10.03", "54.03", "1.45", "0.1454", 5),
    @(3, "005270", "太平改革红利精选灵活配置混合", "1.87", "88.32", "7.74", "0.1447", 3),
    @(4, "011390", "华安添祥6个月持有期混合型证券投资基金", "8.25", "33.54", "1.37", "0.1130", 3),
    @(5, "010896", "太平价值增长股票A", "1.18", "83.63", "7.37", "0.0870", 3),
    @(6, "010897", "太平价值增长股票C", "1.01", "83.63", "7.37", "0.0744", 3),
    @(7, "005695", "华安睿明两年定期开放灵活配置混合A", "1.98", "93.49", "2.60", "0.0515", 9),
    @(8, "005696", "华安睿明两年定期开放灵活配置混合C", "0.10", "93.49", "2.60", "0.0026", 9)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q1new.Range("A$r").Value = $row[0]

    # Columns B-G are stored as text in this workbook (fund codes with
    # leading zeros, decimal figures kept as literal strings) - force text
    # formatting before assigning so Excel doesn't re-interpret them as
    # numbers. (Note: use string concatenation, not "B$r:G$r" - PowerShell
    # parses the colon right after an interpolated variable as a scope
    # qualifier and silently truncates the string.)
    $q1new.Range("B" + $r + ":G" + $r).NumberFormat = "@"
    $q1new.Range("B$r").Value = $row[1]
    $q1new.Range("C$r").Value = $row[2]
    $q1new.Range("D$r").Value = $row[3]
    $q1new.Range("E$r").Value = $row[4]
    $q1new.Range("F$r").Value = $row[5]
    $q1new.Range("G$r").Value = $row[6]

    # Position rank stays numeric.
    $q1new.Range("H$r").Value = $row[7]
}

# --- 2) prepend the 2022-Q1 aggregate row to "总计" --------------------
# No sheets were added/removed since the last fetch above, but re-fetch
# anyway to be explicit/safe.
$tot = $wb.Worksheets.Item("总计")

$tot.Rows.Item(2).Insert()

# The freshly inserted row has no formatting of its own yet - clone it from
# the row right below (the old top data row), which carries the bordered
# index-column style in column A.
$tot.Range("A3:D3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 9
$tot.Range("D2").Value = 1.59

# Renumber the existing rows' index column (they each shift down by one
# position now that 2022-Q1 leads the table).
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
